$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [DateTime]::FromOADate(45243)

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
